# TMTC0032668_VerifyActivityIsLinkedToTheRelatedEngagement
# Add a new "AddContact" worksheet (after "AddOpportunity") with a header
# row and one sample data row, then select it (making it the active tab).

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the last existing sheet ("AddOpportunity").
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "AddContact"

# Header row (bold).
$headers = @("Contact", "Role", "Party", "Type1", "ClientType", "Contact2", "Type2", "HLContact")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $newSheet.Cells.Item(1, $i + 1)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
}

# Sample data row.
$data = @("Chris Lord", "Board of Directors", "Buyer", "External", "Client", "Emma Watson", "Client", "Sonika Goyal")
for ($i = 0; $i -lt $data.Length; $i++) {
    $cell = $newSheet.Cells.Item(2, $i + 1)
    $cell.Value = $data[$i]
}

# Match the author's final selection/active-sheet state.
$newSheet.Range("J14").Select() | Out-Null
